$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty "version002" result cells (C2:C7, C11)
$ws.Range("C2").Value = "SUCCESS - No ODE string stored"
$ws.Range("C3").Value = "SUCCESS - No ODE string stored"
$ws.Range("C4").Value = "SUCCESS - No ODE string stored"
$ws.Range("C5").Value = "FAILED: Error while executing the code: name 'N' is not defined..."
$ws.Range("C6").Value = "FAILED: Error while executing the code: name 'l_a' is not defined..."
$ws.Range("C7").Value = "FAILED: Error while executing the code: can't multiply sequence by non-int of type 'Q'..."
$ws.Range("C11").Value = "SUCCESS - No ODE string stored"

# Replace the stored ODE string in C12 with the new run's result
$ws.Range("C12").Value = "FAILED: Error while executing the code: 'Symbol' object is not callable..."

# Append new rows (13-28) of results for the iterative, 2-step prompting run
$ws.Range("A13").Value = "BIOMD0000000964"
$ws.Range("C13").Value = "SUCCESS - No ODE string stored"

$ws.Range("A14").Value = "BIOMD0000000970"
$ws.Range("C14").Value = "SUCCESS - No ODE string stored"

$ws.Range("A15").Value = "BIOMD0000000974"
$ws.Range("C15").Value = "SUCCESS - No ODE string stored"

$ws.Range("A16").Value = "BIOMD0000000963"
$ws.Range("C16").Value = "SUCCESS - No ODE string stored"

$ws.Range("A17").Value = "BIOMD0000000977"
$ws.Range("C17").Value = "SUCCESS - No ODE string stored"

$ws.Range("A18").Value = "odes_to_mira_SEVITHR"
$ws.Range("C18").Value = "SUCCESS - No ODE string stored"

$ws.Range("A19").Value = "BIOMD0000000984"
$ws.Range("C19").Value = "SUCCESS - No ODE string stored"

$ws.Range("A20").Value = "BIOMD0000000978"
$ws.Range("C20").Value = "SUCCESS - No ODE string stored"

$ws.Range("A21").Value = "Dec_2024_epi_scenario1_modelC"
$ws.Range("C21").Value = "SUCCESS - No ODE string stored"

$ws.Range("A22").Value = "Dec_2024_epi_scenario1_modelB"
$ws.Range("C22").Value = "SUCCESS - No ODE string stored"

$ws.Range("A23").Value = "BIOMD0000000972"
$ws.Range("C23").Value = "FAILED: Error while executing the code: 'Symbol' object is not callable..."

$ws.Range("A24").Value = "BIOMD0000000971"
$ws.Range("C24").Value = "FAILED: Error while executing the code: name 'theta' is not defined..."

$ws.Range("A25").Value = "BIOMD0000000976"
$ws.Range("C25").Value = "FAILED: Error while executing the code: 'Symbol' object is not callable..."

$ws.Range("A26").Value = "BIOMD0000000979"
$ws.Range("C26").Value = "FAILED: Error while executing the code: exec() arg 1 must be a string, bytes or code object..."

$ws.Range("A27").Value = "BIOMD0000000983"
$ws.Range("C27").Value = "FAILED: Error while executing the code: 'Symbol' object is not callable..."

$ws.Range("A28").Value = "Dec_2024_epi_scenario1_modelA"
$ws.Range("C28").Value = "FAILED: Error while executing the code: name 'pi_n' is not defined..."

$wb.Save()
